# Insert a new weekly price record for Jengibre (ginger) at row 91,
# pushing the existing rows 91-101 down to 92-102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 91 (shifts rows 91:101 -> 92:102)
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new record's values.
$ws.Range("A91").Value = 8
$ws.Range("B91").Value = "Terminal La Palmera de La Serena"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44946
$ws.Range("E91").Value = 4
$ws.Range("F91").Value = 100114007
$ws.Range("G91").Value = "Jengibre"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 400
$ws.Range("K91").Value = 13000
$ws.Range("L91").Value = 14000
$ws.Range("M91").Value = 13500
$ws.Range("N91").Value = "`$/caja 13 kilos"
$ws.Range("O91").Value = "Perú"
$ws.Range("P91").Value = 1038
$ws.Range("Q91").Value = 13
$ws.Range("R91").Value = "Hortaliza"
